# Weekly refresh: insert the newest week's two price observations at the
# top of the data block (rows 197-198), pushing all the existing rows
# 197-298 down by two (to 199-300). No other values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 197; this shifts the
# existing rows 197..298 down to 199..300 (and grows the used range to
# A1:R300), exactly mirroring the diff's row-shift pattern.
$ws.Rows("197:198").Insert()

# --- New row 197 ---
$ws.Cells.Item(197, 1).Value  = 7
$ws.Cells.Item(197, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(197, 3).Value  = "Ñuble"
$ws.Cells.Item(197, 4).Value  = 44784
$ws.Cells.Item(197, 5).Value  = 16
$ws.Cells.Item(197, 6).Value  = 100112008
$ws.Cells.Item(197, 7).Value  = "Coliflor"
$ws.Cells.Item(197, 8).Value  = "Sin especificar"
$ws.Cells.Item(197, 9).Value  = "Primera"
$ws.Cells.Item(197, 10).Value = 300
$ws.Cells.Item(197, 11).Value = 900
$ws.Cells.Item(197, 12).Value = 1000
$ws.Cells.Item(197, 13).Value = 950
$ws.Cells.Item(197, 14).Value = "`$/unidad"
$ws.Cells.Item(197, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(197, 16).Value = 950
$ws.Cells.Item(197, 17).Value = 1
$ws.Cells.Item(197, 18).Value = "Hortaliza"

# --- New row 198 ---
$ws.Cells.Item(198, 1).Value  = 7
$ws.Cells.Item(198, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(198, 3).Value  = "Ñuble"
$ws.Cells.Item(198, 4).Value  = 44784
$ws.Cells.Item(198, 5).Value  = 16
$ws.Cells.Item(198, 6).Value  = 100112008
$ws.Cells.Item(198, 7).Value  = "Coliflor"
$ws.Cells.Item(198, 8).Value  = "Sin especificar"
$ws.Cells.Item(198, 9).Value  = "Segunda"
$ws.Cells.Item(198, 10).Value = 200
$ws.Cells.Item(198, 11).Value = 800
$ws.Cells.Item(198, 12).Value = 800
$ws.Cells.Item(198, 13).Value = 800
$ws.Cells.Item(198, 14).Value = "`$/unidad"
$ws.Cells.Item(198, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(198, 16).Value = 800
$ws.Cells.Item(198, 17).Value = 1
$ws.Cells.Item(198, 18).Value = "Hortaliza"

# Give the two new date cells the same date style ("s=2" / numFmt 165)
# used by every other cell in column D.
$ws.Range("D197:D198").NumberFormat = $ws.Range("D199").NumberFormat
